# Electricity Technology Shareweights.xlsx -- apply commit:
# "Updates to electricity exponent and shareweights"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsETS   = $wb.Worksheets.Item("ETS")

# --- Enable iterative calculation (Excel Options > Formulas > Enable iterative
# calculation), matching the workbook's new calcPr iterate settings. ---
$excel.Iteration  = $true
$excel.MaxChange  = 0.00001

# --- ETS sheet: zero out the "hydro" shareweight row (row 5, columns B:AF) ---
$wsETS.Range("B5:AF5").Value = 0

# --- About sheet: document why hydro/crude oil/fuel oil are zeroed out ---
$wsAbout.Range("A20").Value = "We assign zero values to hydro, crude oil and fuel oil, as these plants to not have a pathway"
$wsAbout.Range("A21").Value = "to deployment in the US."

# --- Restore selection/active-cell state seen in the saved workbook ---
[void]$wsETS.Range("B17").Select()
[void]$wsAbout.Range("A22").Select()
